$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep exact text formatting (no numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.412.29'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.844.45'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '240.20'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Value = '0.6276'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = '0.07415'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("D9").Value = '0.2895'
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").Value = '24.73'
$ws.Range("E10").Value = '  +1.15%  '
$ws.Range("D11").Value = '0.07732'
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = '1.843.53'
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("D13").Value = '4.994'
$ws.Range("E13").Value = '  -0.24%  '
$ws.Range("D14").Value = '0.6778'
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = '0.00001009'
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("D16").Value = '82.07'
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("D17").Value = '6.272'
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").Value = '29.471.14'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = '229.10'
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = '12.33'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.18%  '
$ws.Range("D22").Value = '7.453'
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = '1.002'
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").Value = '158.92'
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("D25").Value = '8.472'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '0.1355'
$ws.Range("E26").Value = '  -2.59%  '
$ws.Range("D27").Value = '17.48'
$ws.Range("E27").Value = '  -0.94%  '
$ws.Range("D28").Value = '0.06640'
$ws.Range("E28").Value = '  +17.27%  '
$ws.Range("D29").Value = '1.459'
$ws.Range("E29").Value = '  +1.83%  '
$ws.Range("D30").Value = '1.491'
$ws.Range("E30").Value = '  +0.94%  '
$ws.Range("D31").Value = '4.080'
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("D32").Value = '4.066'
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = '1.841'
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("D34").Value = '1.137'
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("D35").Value = '0.7012'
$ws.Range("E35").Value = '  +0.31%  '
$ws.Range("D36").Value = '2.568'
$ws.Range("E36").Value = '  -0.42%  '
$ws.Range("D37").Value = '0.01855'
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("D38").Value = '2.821'
$ws.Range("E38").Value = '  +3.86%  '
$ws.Range("D39").Value = '1.246.46'
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").Value = '6.743'
$ws.Range("E40").Value = '  +5.12%  '
$ws.Range("D41").Value = '0.9399'
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("D43").Value = '2.024.07'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("D44").Value = '100.93'
$ws.Range("E44").Value = '  -0.42%  '
$ws.Range("D45").Value = '65.93'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("E46").Value = '  +4.61%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.057'
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '1.707'
$ws.Range("E48").Value = '  +2.11%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1151'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '8.965'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("B51").Value = 'TheSandbox'
$ws.Range("C51").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D51").Value = '0.3913'
$ws.Range("E51").Value = '  -0.71%  '
